$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply scraped price / volume updates from the latest cryptos refresh.
# Force text format first so numeric-looking strings (e.g. "567.55",
# "0.123", scientific-looking "0.0000169") are preserved as literal text
# exactly like the source inlineStr cells, instead of being coerced to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.373.56'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.376.96'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.55'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.77'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -6.42%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.378.18'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.01%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.76%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.74%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.02%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.40%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.952.06'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.09%  '
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.123'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.94%  '
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.93'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.374.93'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.71%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000169'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.54%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '60.523.44'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.96%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.78'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -4.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.98'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -5.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '385.98'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.89%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.41%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.96'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.31%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.21%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -8.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.518.92'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.96%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.60%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.36'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -4.87%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.14%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -9.00%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.43'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.406.64'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.85%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '167.99'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.14%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.85'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.94'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -5.26%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.94%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0767'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.33%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '26.99'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.77%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.778'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.83%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '41.33'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.39%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.68'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.65%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.513.49'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.57%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -4.66%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.67%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.27%  '
